$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the diagonal values (second isosceles numbers) per the commit.
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 20

$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 20

$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 20

$ws.Range("H5").Value = 19
$ws.Range("I5").Value = 20

$ws.Range("I6").Value = 19
$ws.Range("J6").Value = 20

$ws.Range("J7").Value = 19
$ws.Range("K7").Value = 20

$ws.Range("N8").Value = 20

$ws.Range("N9").Value = 20

$ws.Range("M10").Value = 16
$ws.Range("N10").Value = 20

$ws.Range("N11").Value = 16
$ws.Range("O11").Value = 20

$ws.Range("N12").Value = 17
$ws.Range("O12").Value = 20

# Update the selection, as recorded in the saved workbook.
$ws.Range("P14").Select()
